$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("I1").Value = "Other found locations"

# Row 2 - update Authors (E2) and add Other found locations (I2)
$ws.Range("E2").Value = "[Luciano%Gattinoni%gattinoniluciano@gmail.com%2,   Davide%Chiumello%NULL%3,   Sandra%Rossi%NULL%0]"
$ws.Range("I2").Value = "_PMC_Springer"

# Row 3 - update Authors (E3) and add Other found locations (I3)
$ws.Range("E3").Value = "[Luciano%Gattinoni%NULL%0,   Silvia%Coppola%NULL%2,   Silvia%Coppola%NULL%0,   Massimo%Cressoni%NULL%1,   Mattia%Busana%NULL%2,   Mattia%Busana%NULL%0,   Sandra%Rossi%NULL%0,   Sandra%Rossi%NULL%0,   Davide%Chiumello%NULL%0,   Davide%Chiumello%NULL%0]"
$ws.Range("I3").Value = "_PMC"

# Row 4 - update Authors (E4) and add Other found locations (I4)
$ws.Range("E4").Value = "[Khai%Tran%NULL%1,   Karen%Cimon%NULL%1,   Melissa%Severn%NULL%1,   Carmem L.%Pessoa-Silva%NULL%1,   John%Conly%NULL%1,   Malcolm Gracie%Semple%NULL%2,   Malcolm Gracie%Semple%NULL%0]"
$ws.Range("I4").Value = "_PMC"
